# Regenerate save_data column G ("K") values for rows 2-24.
# Corresponds to: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" — the recalculated K values replace the old Strike#
# derived values in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKValues = @{
    2  = 0
    3  = 4
    4  = 5
    5  = 7
    6  = 1
    7  = 3
    8  = 4
    9  = 6
    10 = 8
    11 = 2
    12 = 2
    13 = 3
    14 = 3
    15 = 2
    16 = 4
    17 = 3
    18 = 4
    19 = 1
    20 = 1
    21 = 4
    22 = 5
    23 = 2
    24 = 4
}

foreach ($row in $newKValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newKValues[$row]
}
